# Auto-generated edit script: updates cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price cells whose new value would otherwise be
# auto-detected as a number by Excel (losing formatting like trailing zeros).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values.
$ws.Range("D2").Value = "25.786.57"
$ws.Range("E2").Value = "  +5.97%  "
$ws.Range("D3").Value = "1.707.26"
$ws.Range("E3").Value = "  +3.61%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "331.31"
$ws.Range("E5").Value = "  +6.60%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "0.3680"
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("D8").Value = "48.40"
$ws.Range("E8").Value = "  +3.59%  "
$ws.Range("D9").Value = "0.3304"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").Value = "1.170"
$ws.Range("E10").Value = "  +5.12%  "
$ws.Range("D11").Value = "0.07350"
$ws.Range("E11").Value = "  +5.19%  "
$ws.Range("D12").Value = "0.9986"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "6.201"
$ws.Range("E13").Value = "  +4.96%  "
$ws.Range("D14").Value = "19.95"
$ws.Range("E14").Value = "  +3.76%  "
$ws.Range("D15").Value = "6.863"
$ws.Range("E15").Value = "  +4.90%  "
$ws.Range("D16").Value = "1.698.95"
$ws.Range("E16").Value = "  +3.31%  "
$ws.Range("D17").Value = "0.00001073"
$ws.Range("E17").Value = "  +4.21%  "
$ws.Range("D18").Value = "0.06614"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "81.32"
$ws.Range("E19").Value = "  +4.69%  "
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "6.065"
$ws.Range("E21").Value = "  +3.04%  "
$ws.Range("D22").Value = "16.19"
$ws.Range("E22").Value = "  +4.62%  "
$ws.Range("D23").Value = "12.99"
$ws.Range("E23").Value = "  +4.76%  "
$ws.Range("D24").Value = "25.749.77"
$ws.Range("E24").Value = "  +5.77%  "
$ws.Range("D25").Value = "2.465"
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "2.485"
$ws.Range("E26").Value = "  +9.17%  "
$ws.Range("D27").Value = "149.67"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("D28").Value = "19.15"
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").Value = "1.302"
$ws.Range("E29").Value = "  +12.18%  "
$ws.Range("D30").Value = "1.888.28"
$ws.Range("E30").Value = "  +3.31%  "
$ws.Range("D31").Value = "128.03"
$ws.Range("E31").Value = "  +3.73%  "
$ws.Range("D32").Value = "4.112"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "5.956"
$ws.Range("E33").Value = "  +6.57%  "
$ws.Range("D34").Value = "0.08494"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("D35").Value = "1.675"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "12.88"
$ws.Range("E36").Value = "  +7.43%  "
$ws.Range("D37").Value = "5.321"
$ws.Range("E37").Value = "  +3.56%  "
$ws.Range("D38").Value = "1.272"
$ws.Range("E38").Value = "  +3.09%  "
$ws.Range("D39").Value = "0.06218"
$ws.Range("E39").Value = "  +4.66%  "
$ws.Range("D40").Value = "8.546"
$ws.Range("E40").Value = "  +6.24%  "
$ws.Range("D41").Value = "0.2125"
$ws.Range("E41").Value = "  +4.26%  "
$ws.Range("D42").Value = "0.02259"
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("D43").Value = "14.43"
$ws.Range("E43").Value = "  +16.21%  "
$ws.Range("D44").Value = "0.6119"
$ws.Range("E44").Value = "  +4.95%  "
$ws.Range("D45").Value = "0.9985"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "3.846"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").Value = "0.5839"
$ws.Range("E47").Value = "  +5.18%  "
$ws.Range("D48").Value = "126.35"
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("D49").Value = "2.004"
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("D50").Value = "0.07218"
$ws.Range("E50").Value = "  +5.04%  "
$ws.Range("E51").Value = "  +2.93%  "
